# Auto-generated Excel COM-interop script
# Restructures the workbook: adds 'Player Info' and 'ODI Batting Extra' sheets,
# renames MATCH_CARD_LINK columns to MATCH_CODE and stores bare match codes.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-HeaderCell($cell, $val) {
    $cell.Value = $val
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$battingSheet = $wb.Worksheets.Item('ODI Batting')
$bowlingSheet = $wb.Worksheets.Item('ODI Bowling')

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = 'Player Info'

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = 'ODI Batting Extra'

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE -------------------------
Set-HeaderCell $battingSheet.Cells.Item(1, 4) 'MATCH_CODE'
Set-TextCell $battingSheet.Cells.Item(2, 4) '4017'
Set-TextCell $battingSheet.Cells.Item(3, 4) '4018'
Set-TextCell $battingSheet.Cells.Item(4, 4) '4019'
Set-TextCell $battingSheet.Cells.Item(5, 4) '4034'
Set-TextCell $battingSheet.Cells.Item(6, 4) '4037'
Set-TextCell $battingSheet.Cells.Item(7, 4) '4045'
Set-TextCell $battingSheet.Cells.Item(8, 4) '4050'
Set-TextCell $battingSheet.Cells.Item(9, 4) '4079'
Set-TextCell $battingSheet.Cells.Item(10, 4) '4081'
Set-TextCell $battingSheet.Cells.Item(11, 4) '4082'
Set-TextCell $battingSheet.Cells.Item(12, 4) '4084'
Set-TextCell $battingSheet.Cells.Item(13, 4) '4087'
Set-TextCell $battingSheet.Cells.Item(14, 4) '4103'
Set-TextCell $battingSheet.Cells.Item(15, 4) '4104'
Set-TextCell $battingSheet.Cells.Item(16, 4) '4105'
Set-TextCell $battingSheet.Cells.Item(17, 4) '4110'
Set-TextCell $battingSheet.Cells.Item(18, 4) '4114'
Set-TextCell $battingSheet.Cells.Item(19, 4) '4172'
Set-TextCell $battingSheet.Cells.Item(20, 4) '4174'
Set-TextCell $battingSheet.Cells.Item(21, 4) '4176'
Set-TextCell $battingSheet.Cells.Item(22, 4) '4177'
Set-TextCell $battingSheet.Cells.Item(23, 4) '4178'
Set-TextCell $battingSheet.Cells.Item(24, 4) '4194'
Set-TextCell $battingSheet.Cells.Item(25, 4) '4197'
Set-TextCell $battingSheet.Cells.Item(26, 4) '4201'
Set-TextCell $battingSheet.Cells.Item(27, 4) '4204'
Set-TextCell $battingSheet.Cells.Item(28, 4) '4223'
Set-TextCell $battingSheet.Cells.Item(29, 4) '4225'
Set-TextCell $battingSheet.Cells.Item(30, 4) '4227'
Set-TextCell $battingSheet.Cells.Item(31, 4) '4237'
Set-TextCell $battingSheet.Cells.Item(32, 4) '4238'
Set-TextCell $battingSheet.Cells.Item(33, 4) '4241'
Set-TextCell $battingSheet.Cells.Item(34, 4) '4244'
Set-TextCell $battingSheet.Cells.Item(35, 4) '4247'
Set-TextCell $battingSheet.Cells.Item(36, 4) '4304'
Set-TextCell $battingSheet.Cells.Item(37, 4) '4308'
Set-TextCell $battingSheet.Cells.Item(38, 4) '4324'
Set-TextCell $battingSheet.Cells.Item(39, 4) '4334'
Set-TextCell $battingSheet.Cells.Item(40, 4) '4337'
Set-TextCell $battingSheet.Cells.Item(41, 4) '4340'
Set-TextCell $battingSheet.Cells.Item(42, 4) '4349'
Set-TextCell $battingSheet.Cells.Item(43, 4) '4375'
Set-TextCell $battingSheet.Cells.Item(44, 4) '4376'
Set-TextCell $battingSheet.Cells.Item(45, 4) '4458'
Set-TextCell $battingSheet.Cells.Item(46, 4) '4459'
Set-TextCell $battingSheet.Cells.Item(47, 4) '4472'
Set-TextCell $battingSheet.Cells.Item(48, 4) '4473'
Set-TextCell $battingSheet.Cells.Item(49, 4) '4476'
Set-TextCell $battingSheet.Cells.Item(50, 4) '4586'
Set-TextCell $battingSheet.Cells.Item(51, 4) '4590'
Set-TextCell $battingSheet.Cells.Item(52, 4) '4592'
Set-TextCell $battingSheet.Cells.Item(53, 4) '4634'
Set-TextCell $battingSheet.Cells.Item(54, 4) '4638'

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE -------------------------
Set-HeaderCell $bowlingSheet.Cells.Item(1, 2) 'MATCH_CODE'
Set-TextCell $bowlingSheet.Cells.Item(2, 2) '4017'
Set-TextCell $bowlingSheet.Cells.Item(3, 2) '4018'
Set-TextCell $bowlingSheet.Cells.Item(4, 2) '4019'
Set-TextCell $bowlingSheet.Cells.Item(5, 2) '4034'
Set-TextCell $bowlingSheet.Cells.Item(6, 2) '4037'
Set-TextCell $bowlingSheet.Cells.Item(7, 2) '4045'
Set-TextCell $bowlingSheet.Cells.Item(8, 2) '4050'
Set-TextCell $bowlingSheet.Cells.Item(9, 2) '4079'
Set-TextCell $bowlingSheet.Cells.Item(10, 2) '4081'
Set-TextCell $bowlingSheet.Cells.Item(11, 2) '4082'
Set-TextCell $bowlingSheet.Cells.Item(12, 2) '4084'
Set-TextCell $bowlingSheet.Cells.Item(13, 2) '4087'
Set-TextCell $bowlingSheet.Cells.Item(14, 2) '4103'
Set-TextCell $bowlingSheet.Cells.Item(15, 2) '4104'
Set-TextCell $bowlingSheet.Cells.Item(16, 2) '4105'
Set-TextCell $bowlingSheet.Cells.Item(17, 2) '4110'
Set-TextCell $bowlingSheet.Cells.Item(18, 2) '4114'
Set-TextCell $bowlingSheet.Cells.Item(19, 2) '4172'
Set-TextCell $bowlingSheet.Cells.Item(20, 2) '4174'
Set-TextCell $bowlingSheet.Cells.Item(21, 2) '4176'
Set-TextCell $bowlingSheet.Cells.Item(22, 2) '4177'
Set-TextCell $bowlingSheet.Cells.Item(23, 2) '4178'
Set-TextCell $bowlingSheet.Cells.Item(24, 2) '4194'
Set-TextCell $bowlingSheet.Cells.Item(25, 2) '4197'
Set-TextCell $bowlingSheet.Cells.Item(26, 2) '4201'
Set-TextCell $bowlingSheet.Cells.Item(27, 2) '4204'
Set-TextCell $bowlingSheet.Cells.Item(28, 2) '4223'
Set-TextCell $bowlingSheet.Cells.Item(29, 2) '4225'
Set-TextCell $bowlingSheet.Cells.Item(30, 2) '4237'
Set-TextCell $bowlingSheet.Cells.Item(31, 2) '4238'
Set-TextCell $bowlingSheet.Cells.Item(32, 2) '4241'
Set-TextCell $bowlingSheet.Cells.Item(33, 2) '4244'
Set-TextCell $bowlingSheet.Cells.Item(34, 2) '4247'
Set-TextCell $bowlingSheet.Cells.Item(35, 2) '4308'
Set-TextCell $bowlingSheet.Cells.Item(36, 2) '4324'
Set-TextCell $bowlingSheet.Cells.Item(37, 2) '4334'
Set-TextCell $bowlingSheet.Cells.Item(38, 2) '4337'
Set-TextCell $bowlingSheet.Cells.Item(39, 2) '4340'
Set-TextCell $bowlingSheet.Cells.Item(40, 2) '4349'
Set-TextCell $bowlingSheet.Cells.Item(41, 2) '4375'
Set-TextCell $bowlingSheet.Cells.Item(42, 2) '4376'
Set-TextCell $bowlingSheet.Cells.Item(43, 2) '4458'
Set-TextCell $bowlingSheet.Cells.Item(44, 2) '4459'
Set-TextCell $bowlingSheet.Cells.Item(45, 2) '4472'
Set-TextCell $bowlingSheet.Cells.Item(46, 2) '4473'
Set-TextCell $bowlingSheet.Cells.Item(47, 2) '4476'
Set-TextCell $bowlingSheet.Cells.Item(48, 2) '4586'
Set-TextCell $bowlingSheet.Cells.Item(49, 2) '4590'
Set-TextCell $bowlingSheet.Cells.Item(50, 2) '4592'
Set-TextCell $bowlingSheet.Cells.Item(51, 2) '4634'
Set-TextCell $bowlingSheet.Cells.Item(52, 2) '4638'

# --- Player Info ---------------------------------------------------------
Set-HeaderCell $playerInfo.Cells.Item(1, 1) 'ID'
Set-HeaderCell $playerInfo.Cells.Item(1, 2) 'NAME'
Set-HeaderCell $playerInfo.Cells.Item(1, 3) 'BATTING_HAND'
Set-HeaderCell $playerInfo.Cells.Item(1, 4) 'BOWL_STYLE'
Set-TextCell $playerInfo.Cells.Item(2, 1) '4653'
Set-TextCell $playerInfo.Cells.Item(2, 2) 'Shadab Khan'
Set-TextCell $playerInfo.Cells.Item(2, 3) 'Right Handed'
Set-TextCell $playerInfo.Cells.Item(2, 4) 'Right Arm Leg Break'

# --- ODI Batting Extra ---------------------------------------------------
Set-HeaderCell $battingExtra.Cells.Item(1, 1) 'MATCH_CODE'
Set-HeaderCell $battingExtra.Cells.Item(1, 2) 'BATTING_POSITION'
Set-HeaderCell $battingExtra.Cells.Item(1, 3) 'NUM_4'
Set-HeaderCell $battingExtra.Cells.Item(1, 4) 'NUM_6'
Set-HeaderCell $battingExtra.Cells.Item(1, 5) 'PERCENT_RUNS_OF_TOTAL'
Set-HeaderCell $battingExtra.Cells.Item(1, 6) 'MAN_OF_MATCH'
Set-TextCell $battingExtra.Cells.Item(2, 1) '4247'
Set-TextCell $battingExtra.Cells.Item(2, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(3, 1) '4304'
$battingExtra.Cells.Item(3, 2).Value = 8
Set-TextCell $battingExtra.Cells.Item(3, 3) '0'
Set-TextCell $battingExtra.Cells.Item(3, 4) '0'
Set-TextCell $battingExtra.Cells.Item(3, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(4, 1) '4308'
Set-TextCell $battingExtra.Cells.Item(4, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(5, 1) '4324'
$battingExtra.Cells.Item(5, 2).Value = 8
Set-TextCell $battingExtra.Cells.Item(5, 3) '1'
Set-TextCell $battingExtra.Cells.Item(5, 4) '0'
Set-TextCell $battingExtra.Cells.Item(5, 5) '9.43%'
Set-TextCell $battingExtra.Cells.Item(5, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(6, 1) '4334'
$battingExtra.Cells.Item(6, 2).Value = 9
Set-TextCell $battingExtra.Cells.Item(6, 3) '0'
Set-TextCell $battingExtra.Cells.Item(6, 4) '0'
Set-TextCell $battingExtra.Cells.Item(6, 5) '0.32%'
Set-TextCell $battingExtra.Cells.Item(6, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(7, 1) '4337'
Set-TextCell $battingExtra.Cells.Item(7, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(8, 1) '4340'
$battingExtra.Cells.Item(8, 2).Value = 8
Set-TextCell $battingExtra.Cells.Item(8, 3) '1'
Set-TextCell $battingExtra.Cells.Item(8, 4) '0'
Set-TextCell $battingExtra.Cells.Item(8, 5) '4.78%'
Set-TextCell $battingExtra.Cells.Item(8, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(9, 1) '4349'
$battingExtra.Cells.Item(9, 2).Value = 9
Set-TextCell $battingExtra.Cells.Item(9, 3) '0'
Set-TextCell $battingExtra.Cells.Item(9, 4) '0'
Set-TextCell $battingExtra.Cells.Item(9, 5) '0.32%'
Set-TextCell $battingExtra.Cells.Item(9, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(10, 1) '4375'
$battingExtra.Cells.Item(10, 2).Value = 9
Set-TextCell $battingExtra.Cells.Item(10, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(11, 1) '4376'
$battingExtra.Cells.Item(11, 2).Value = 9
Set-TextCell $battingExtra.Cells.Item(11, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(12, 1) '4458'
$battingExtra.Cells.Item(12, 2).Value = 7
Set-TextCell $battingExtra.Cells.Item(12, 3) '2'
Set-TextCell $battingExtra.Cells.Item(12, 4) '1'
Set-TextCell $battingExtra.Cells.Item(12, 5) '12.04%'
Set-TextCell $battingExtra.Cells.Item(12, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(13, 1) '4459'
$battingExtra.Cells.Item(13, 2).Value = 6
Set-TextCell $battingExtra.Cells.Item(13, 3) '1'
Set-TextCell $battingExtra.Cells.Item(13, 4) '0'
Set-TextCell $battingExtra.Cells.Item(13, 5) '4.01%'
Set-TextCell $battingExtra.Cells.Item(13, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(14, 1) '4472'
$battingExtra.Cells.Item(14, 2).Value = 7
Set-TextCell $battingExtra.Cells.Item(14, 3) '0'
Set-TextCell $battingExtra.Cells.Item(14, 4) '1'
Set-TextCell $battingExtra.Cells.Item(14, 5) '21.28%'
Set-TextCell $battingExtra.Cells.Item(14, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(15, 1) '4473'
$battingExtra.Cells.Item(15, 2).Value = 7
Set-TextCell $battingExtra.Cells.Item(15, 3) '3'
Set-TextCell $battingExtra.Cells.Item(15, 4) '0'
Set-TextCell $battingExtra.Cells.Item(15, 5) '10.77%'
Set-TextCell $battingExtra.Cells.Item(15, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(16, 1) '4476'
$battingExtra.Cells.Item(16, 2).Value = 8
Set-TextCell $battingExtra.Cells.Item(16, 3) '0'
Set-TextCell $battingExtra.Cells.Item(16, 4) '0'
Set-TextCell $battingExtra.Cells.Item(16, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(17, 1) '4586'
Set-TextCell $battingExtra.Cells.Item(17, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(18, 1) '4590'
Set-TextCell $battingExtra.Cells.Item(18, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(19, 1) '4592'
$battingExtra.Cells.Item(19, 2).Value = 7
Set-TextCell $battingExtra.Cells.Item(19, 3) '4'
Set-TextCell $battingExtra.Cells.Item(19, 4) '3'
Set-TextCell $battingExtra.Cells.Item(19, 5) '31.97%'
Set-TextCell $battingExtra.Cells.Item(19, 6) 'YES'
Set-TextCell $battingExtra.Cells.Item(20, 1) '4634'
Set-TextCell $battingExtra.Cells.Item(20, 6) 'NO'
Set-TextCell $battingExtra.Cells.Item(21, 1) '4638'
$battingExtra.Cells.Item(21, 2).Value = 6
Set-TextCell $battingExtra.Cells.Item(21, 6) 'NO'

Write-Output 'Edit complete'
